$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("G5").Value = 2.05
$ws.Range("I5").Value = 4.2
$ws.Range("L5").Value = 4.75
$ws.Range("M5").Value = 1.11
$ws.Range("N5").Value = 6.5
$ws.Range("AD5").Value = 17
$ws.Range("AE5").Value = 19
$ws.Range("AI5").Value = 19

# Row 8
$ws.Range("M8").Value = 1.14
$ws.Range("N8").Value = 5.5

# Row 9
$ws.Range("G9").Value = 1.91
$ws.Range("H9").Value = 3.25
$ws.Range("J9").Value = 2.75
$ws.Range("M9").Value = 1.11
$ws.Range("N9").Value = 6.5

# Row 10
$ws.Range("N10").Value = 8
$ws.Range("O10").Value = 1.4
$ws.Range("P10").Value = 2.75

# Row 15
$ws.Range("L15").Value = 3.5

# Row 19
$ws.Range("G19").Value = 3.45
$ws.Range("H19").Value = 3.3
$ws.Range("I19").Value = 2
$ws.Range("J19").Value = 3.95
$ws.Range("K19").Value = 2.07
$ws.Range("L19").Value = 2.62
$ws.Range("O19").Value = 1.26
$ws.Range("P19").Value = 3.15
$ws.Range("Q19").Value = 1.78
$ws.Range("R19").Value = 1.82
$ws.Range("U19").Value = 2.8
$ws.Range("V19").Value = 1.33
$ws.Range("W19").Value = 1.4
$ws.Range("X19").Value = 2.55
$ws.Range("Y19").Value = 1.65
$ws.Range("Z19").Value = 1.98
$ws.Range("AA19").Value = 11
$ws.Range("AB19").Value = 19.5
$ws.Range("AD19").Value = 50
$ws.Range("AE19").Value = 30
$ws.Range("AF19").Value = 35
$ws.Range("AG19").Value = 10.25
$ws.Range("AH19").Value = 6.5
$ws.Range("AI19").Value = 13
$ws.Range("AJ19").Value = 55
$ws.Range("AK19").Value = 400
$ws.Range("AL19").Value = 7.9
$ws.Range("AM19").Value = 10
$ws.Range("AP19").Value = 15.5
$ws.Range("AQ19").Value = 25

# Row 22
$ws.Range("G22").Value = 2
$ws.Range("I22").Value = 3.4
$ws.Range("L22").Value = 3.75
$ws.Range("S22").Value = 1.93
$ws.Range("T22").Value = 1.93
$ws.Range("AO22").Value = 41

# Row 23
$ws.Range("S23").Value = 1.83
$ws.Range("T23").Value = 2.03
$ws.Range("U23").Value = 2.2
$ws.Range("V23").Value = 1.62

# Row 24
$ws.Range("U24").Value = 5
$ws.Range("V24").Value = 1.17
$ws.Range("AR24").Value = 1.95
$ws.Range("AS24").Value = 1.9

# Row 27
$ws.Range("AB27").Value = 9.5
$ws.Range("AC27").Value = 8
$ws.Range("AI27").Value = 13.5
$ws.Range("AK27").Value = 400
